$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column K: "Work Type" header + "Brickwork (Masonry)" value
$ws.Range("K1").Value = "Work Type"
$ws.Range("K2").Value = "Brickwork (Masonry)"

# Match header formatting used for the other header cells (bold, bordered, centered)
$ws.Range("K1").Font.Bold = $true
$ws.Range("K1").Font.Name = "Calibri"
$ws.Range("K1").Font.Size = 11
$ws.Range("K1").Borders.LineStyle = 1
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").VerticalAlignment = -4160

$ws.Range("K2").Font.Name = "Calibri"

$ws.Columns.Item(11).ColumnWidth = 17.1666666667

$ws.Range("N5").Select()
